$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Update the report date (shared across all date cells E2:E7 on both sheets)
$ws1.Range("E2:E7").Value = "2023-08-25"
$ws2.Range("E2:E7").Value = "2023-08-25"

# checkReport sheet: selection moves from E2:E7 to E3:E7, and it's no longer the active tab
$ws1.Range("E3:E7").Select()

# viewMore sheet becomes the active tab, with its selection anchored near E7/E12
$ws2.Activate()
$ws2.Range("E12").Select()
